# Auto-generated edit script: update crypto price/volume table
# Applies the diff: updates Price (D) and Volume(1h) (E) columns for rows 2-51,
# and for rows 29/30 and 42/43 also updates Coin (B) and Link (C) since those rows were reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.048.19"
$ws.Range("E2").Value = "'  -2.97%  "
$ws.Range("D3").Value = "'1.903.08"
$ws.Range("E3").Value = "'  -3.55%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  -0.99%  "
$ws.Range("D5").Value = "'327.05"
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("E6").Value = "'  -0.95%  "
$ws.Range("D7").Value = "'0.4609"
$ws.Range("E7").Value = "'  -4.66%  "
$ws.Range("D8").Value = "'0.3970"
$ws.Range("E8").Value = "'  -2.94%  "
$ws.Range("D9").Value = "'51.87"
$ws.Range("E9").Value = "'  -4.08%  "
$ws.Range("D10").Value = "'0.08326"
$ws.Range("E10").Value = "'  -4.11%  "
$ws.Range("D11").Value = "'1.039"
$ws.Range("E11").Value = "'  -3.13%  "
$ws.Range("D12").Value = "'21.82"
$ws.Range("E12").Value = "'  -2.98%  "
$ws.Range("D13").Value = "'1.912.61"
$ws.Range("E13").Value = "'  -1.37%  "
$ws.Range("D14").Value = "'7.361"
$ws.Range("E14").Value = "'  -5.41%  "
$ws.Range("D15").Value = "'6.018"
$ws.Range("E15").Value = "'  -4.48%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "'  -1.21%  "
$ws.Range("D17").Value = "'89.07"
$ws.Range("E17").Value = "'  -2.29%  "
$ws.Range("D18").Value = "'0.00001064"
$ws.Range("E18").Value = "'  -1.14%  "
$ws.Range("D19").Value = "'0.06601"
$ws.Range("E19").Value = "'  -0.41%  "
$ws.Range("D20").Value = "'17.77"
$ws.Range("E20").Value = "'  -6.10%  "
$ws.Range("E21").Value = "'  -1.09%  "
$ws.Range("D22").Value = "'5.699"
$ws.Range("E22").Value = "'  -2.09%  "
$ws.Range("D23").Value = "'28.053.84"
$ws.Range("E23").Value = "'  -3.03%  "
$ws.Range("D24").Value = "'11.14"
$ws.Range("E24").Value = "'  -4.38%  "
$ws.Range("D25").Value = "'2.312"
$ws.Range("E25").Value = "'  +0.86%  "
$ws.Range("D26").Value = "'2.137.81"
$ws.Range("E26").Value = "'  -1.37%  "
$ws.Range("D27").Value = "'153.41"
$ws.Range("E27").Value = "'  -0.89%  "
$ws.Range("D28").Value = "'19.95"
$ws.Range("E28").Value = "'  -2.37%  "
$ws.Range("B29").Value = "'LidoDAOToken"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.118"
$ws.Range("E29").Value = "'  -3.05%  "
$ws.Range("B30").Value = "'InternetComputer(DFINITY)"
$ws.Range("C30").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'5.708"
$ws.Range("E30").Value = "'  -6.27%  "
$ws.Range("D31").Value = "'123.38"
$ws.Range("E31").Value = "'  -1.63%  "
$ws.Range("D32").Value = "'0.9686"
$ws.Range("E32").Value = "'  -4.75%  "
$ws.Range("D33").Value = "'0.09583"
$ws.Range("E33").Value = "'  -1.19%  "
$ws.Range("D34").Value = "'1.464"
$ws.Range("E34").Value = "'  -1.14%  "
$ws.Range("D35").Value = "'3.621"
$ws.Range("E35").Value = "'  -2.26%  "
$ws.Range("D36").Value = "'5.497"
$ws.Range("E36").Value = "'  -4.19%  "
$ws.Range("D37").Value = "'1.265"
$ws.Range("E37").Value = "'  -1.83%  "
$ws.Range("D38").Value = "'0.02278"
$ws.Range("E38").Value = "'  -3.79%  "
$ws.Range("D39").Value = "'8.675"
$ws.Range("E39").Value = "'  -2.70%  "
$ws.Range("D40").Value = "'0.06130"
$ws.Range("E40").Value = "'  -2.56%  "
$ws.Range("D41").Value = "'0.6117"
$ws.Range("E41").Value = "'  -3.02%  "
$ws.Range("B42").Value = "'Frax"
$ws.Range("C42").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "'  -1.00%  "
$ws.Range("B43").Value = "'Aptos"
$ws.Range("C43").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'10.81"
$ws.Range("E43").Value = "'  -3.92%  "
$ws.Range("D44").Value = "'0.1898"
$ws.Range("E44").Value = "'  -1.80%  "
$ws.Range("D45").Value = "'1.305"
$ws.Range("E45").Value = "'  -3.01%  "
$ws.Range("D46").Value = "'0.5831"
$ws.Range("E46").Value = "'  -3.30%  "
$ws.Range("D47").Value = "'12.78"
$ws.Range("E47").Value = "'  -2.38%  "
$ws.Range("D48").Value = "'1.998"
$ws.Range("E48").Value = "'  -5.17%  "
$ws.Range("D49").Value = "'3.433"
$ws.Range("E49").Value = "'  -0.63%  "
$ws.Range("D50").Value = "'0.06906"
$ws.Range("E50").Value = "'  +0.42%  "
$ws.Range("D51").Value = "'110.32"
$ws.Range("E51").Value = "'  -1.00%  "
